$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers (e.g. "578.63")
# must be force-typed as Text (matching the source workbook, where every
# Price cell is stored as text) by setting NumberFormat to "@" before the
# assignment - otherwise Excel auto-converts them to numbers and silently
# drops things like trailing zeros (e.g. "1.00" -> 1).
$textForceCells = @("D5","D6","D7","D8","D11","D14","D20","D24","D25","D26","D29","D31","D32","D38","D39","D43","D48","D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.775.40'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '3.320.36'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '578.63'
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').Value = '175.01'
$ws.Range('E6').Value = '  -4.34%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').Value = '3.317.33'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').Value = '  -0.51%  '
$ws.Range('D11').Value = '0.574'
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').Value = '660.37'
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('D15').Value = '3.863.31'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').Value = '67.645.00'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('D19').Value = '3.327.00'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '17.38'
$ws.Range('E20').Value = '  -1.70%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('E23').Value = '  +5.12%  '
$ws.Range('D24').Value = '16.88'
$ws.Range('E24').Value = '  -4.29%  '
$ws.Range('D25').Value = '98.13'
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('D26').Value = '3.84'
$ws.Range('E26').Value = '  -3.89%  '
$ws.Range('E27').Value = '  -4.09%  '
$ws.Range('E28').Value = '  -3.56%  '
$ws.Range('D29').Value = '33.16'
$ws.Range('E29').Value = '  +2.33%  '
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('D31').Value = '7.22'
$ws.Range('E31').Value = '  +7.88%  '
$ws.Range('D32').Value = '567.03'
$ws.Range('E32').Value = '  -6.30%  '
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '3.666.41'
$ws.Range('E36').Value = '  -7.11%  '
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('D38').Value = '3.24'
$ws.Range('E38').Value = '  -7.71%  '
$ws.Range('D39').Value = '34.44'
$ws.Range('E39').Value = '  +5.54%  '
$ws.Range('E40').Value = '  +1.05%  '
$ws.Range('E41').Value = '  -2.23%  '
$ws.Range('E42').Value = '  -4.86%  '
$ws.Range('D43').Value = '3.32'
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('E44').Value = '  -1.77%  '
$ws.Range('D45').Value = '0.0₃0659'
$ws.Range('E45').Value = '  -3.98%  '
$ws.Range('E46').Value = '  -2.72%  '
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('D48').Value = '0.128'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('E50').Value = '  -3.26%  '
$ws.Range('D51').Value = '127.68'
$ws.Range('E51').Value = '  -2.53%  '
